$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New text values must be entered in this precise order so that the
#     sharedStrings table gets the same new-string ordering as the source
#     workbook: Varun, Shrikanth, Sameeksha, banglore, snaik2187@..,
#     Shrikant.Salke@.., vsrinivasan2203@.. ---
$ws.Range("B105").Value = "Varun"
$ws.Range("B106").Value = "Shrikanth"
$ws.Range("B107").Value = "Sameeksha"
$ws.Range("E105").Value = "banglore"
$ws.Range("E106").Value = "banglore"
$ws.Range("E107").Value = "banglore"
$ws.Range("G107").Value = "snaik2187@altimetrik.com"
$ws.Range("G106").Value = "Shrikant.Salke@altimetrik.com"
$ws.Range("G105").Value = "vsrinivasan2203@altimetrik.com"

# --- Remaining cell values (reuse existing shared strings / numbers) ---
$ws.Range("A105").Value = 2203
$ws.Range("C105").Value = "ETV"
$ws.Range("D105").Value = "Intuit"
$ws.Range("F105").Value = 917410186388
$ws.Range("F105").NumberFormat = "0"
$ws.Range("H105").Value = "Ramu"
$ws.Range("I105").Value = "ramu@yahoo.com"
$ws.Range("J105").Value = 2399
$ws.Range("K105").Value = 919884899055
$ws.Range("K105").NumberFormat = "0"

$ws.Range("A106").Value = 3165
$ws.Range("C106").Value = "ETV"
$ws.Range("D106").Value = "Intuit"
$ws.Range("F106").Value = 919448493932
$ws.Range("F106").NumberFormat = "0"
$ws.Range("H106").Value = "Ramu"
$ws.Range("I106").Value = "ramu@yahoo.com"
$ws.Range("J106").Value = 2399
$ws.Range("K106").Value = 919884899055
$ws.Range("K106").NumberFormat = "0"

$ws.Range("A107").Value = 2187
$ws.Range("C107").Value = "ETV"
$ws.Range("D107").Value = "Intuit"
$ws.Range("F107").Value = 919611666426
$ws.Range("F107").NumberFormat = "0"
$ws.Range("H107").Value = "Ramu"
$ws.Range("I107").Value = "ramu@yahoo.com"
$ws.Range("J107").Value = 2399
$ws.Range("K107").Value = 919884899055
$ws.Range("K107").NumberFormat = "0"

# --- Hyperlinks (actual OS-level hyperlinks); re-apply the "Hyperlink"
#     cell style afterwards since Hyperlinks.Add touches formatting ---
$ws.Hyperlinks.Add($ws.Range("G107"), "mailto:snaik2187@altimetrik.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I105"), "mailto:ramu@yahoo.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I106"), "mailto:ramu@yahoo.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I107"), "mailto:ramu@yahoo.com") | Out-Null

$ws.Range("G107").Style = "Hyperlink"
$ws.Range("I105").Style = "Hyperlink"
$ws.Range("I106").Style = "Hyperlink"
$ws.Range("I107").Style = "Hyperlink"

# --- Selection / view state ---
$ws.Range("A105:H107").Select() | Out-Null
